$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.275.45'
$ws.Range('E2').Value = '  +4.93%  '
$ws.Range('D3').Value = '2.280.90'
$ws.Range('E3').Value = '  +4.55%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'255.54"
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = "'0.643"
$ws.Range('E6').Value = '  +2.71%  '
$ws.Range('D7').Value = "'72.78"
$ws.Range('E7').Value = '  +7.16%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = "'0.683"
$ws.Range('E8').Value = '  +19.76%  '
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').Value = "'1.00"
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = "'40.77"
$ws.Range('D11').Value = "'0.0984"
$ws.Range('E11').Value = '  +5.67%  '
$ws.Range('D12').Value = "'59.42"
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').Value = "'7.48"
$ws.Range('E13').Value = '  +6.82%  '
$ws.Range('D14').Value = "'0.104"
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = '2.620.55'
$ws.Range('E15').Value = '  +4.50%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = "'0.898"
$ws.Range('E16').Value = '  +3.92%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = "'15.01"
$ws.Range('E17').Value = '  +4.66%  '
$ws.Range('D18').Value = '2.279.52'
$ws.Range('E18').Value = '  +3.84%  '
$ws.Range('D19').Value = '43.216.74'
$ws.Range('E19').Value = '  +4.97%  '
$ws.Range('D20').Value = '0.0₃0988'
$ws.Range('E20').Value = '  +3.64%  '
$ws.Range('D21').Value = "'6.31"
$ws.Range('E21').Value = '  +2.78%  '
$ws.Range('D22').Value = "'73.95"
$ws.Range('E22').Value = '  +2.87%  '
$ws.Range('D23').Value = "'237.99"
$ws.Range('E23').Value = '  +2.54%  '
$ws.Range('D24').Value = "'2.13"
$ws.Range('E24').Value = '  +5.68%  '
$ws.Range('D25').Value = "'3.93"
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('D26').Value = "'11.87"
$ws.Range('E26').Value = '  +1.61%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = "'2.49"
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').Value = "'3.71"
$ws.Range('E29').Value = '  +0.61%  '
$ws.Range('E30').Value = '  +2.40%  '
$ws.Range('D31').Value = "'168.44"
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').Value = "'21.38"
$ws.Range('E32').Value = '  +3.88%  '
$ws.Range('E33').Value = '  +10.90%  '
$ws.Range('D34').Value = "'6.17"
$ws.Range('E34').Value = '  +13.30%  '
$ws.Range('D35').Value = "'0.0789"
$ws.Range('E35').Value = '  +5.69%  '
$ws.Range('E36').Value = '  +3.01%  '
$ws.Range('D37').Value = "'29.16"
$ws.Range('E37').Value = '  +10.75%  '
$ws.Range('D38').Value = "'4.76"
$ws.Range('E38').Value = '  +3.98%  '
$ws.Range('D39').Value = "'4.19"
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('D40').Value = "'0.0323"
$ws.Range('E40').Value = '  +9.25%  '
$ws.Range('E41').Value = '  +5.57%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').Value = "'5.96"
$ws.Range('E42').Value = '  +5.73%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').Value = "'12.60"
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('D44').Value = "'64.86"
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('D45').Value = "'5.00"
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('E46').Value = '  +3.04%  '
$ws.Range('D47').Value = "'9.05"
$ws.Range('E47').Value = '  +5.55%  '
$ws.Range('D48').Value = "'0.104"
$ws.Range('E48').Value = '  +3.10%  '
$ws.Range('D49').Value = "'1.21"
$ws.Range('E49').Value = '  +2.07%  '
$ws.Range('D50').Value = "'1.01"
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = "'4.44"
$ws.Range('E51').Value = '  +5.44%  '
